$wb = $excel.ActiveWorkbook

# Scheduled market-data refresh: update price/profit columns (H-N) for the
# affected leve rows across all job sheets, per the latest Universalis pull.

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1190.2424
$ws.Range("I98").Value = 942.7143
$ws.Range("K98").Value = 942.7143
$ws.Range("M98").Value = 555.2857
# Row 122
$ws.Range("H122").Value = 1190.2424
$ws.Range("I122").Value = 942.7143
$ws.Range("K122").Value = 2828.1429
$ws.Range("M122").Value = -378.1428999999998
# Row 135
$ws.Range("H135").Value = 3639.6155
$ws.Range("I135").Value = 963.64703
$ws.Range("K135").Value = 8672.823269999999
$ws.Range("M135").Value = -6137.823269999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3352.96
$ws.Range("I32").Value = 2878.0833
$ws.Range("J32").Value = 14750
$ws.Range("K32").Value = 2878.0833
$ws.Range("L32").Value = 14750
$ws.Range("M32").Value = -2591.0833
$ws.Range("N32").Value = -15324
# Row 45
$ws.Range("H45").Value = 1102.591
$ws.Range("I45").Value = 1067.9
$ws.Range("J45").Value = 1449.5
$ws.Range("K45").Value = 1067.9
$ws.Range("L45").Value = 1449.5
$ws.Range("M45").Value = -690.9000000000001
$ws.Range("N45").Value = -2203.5
# Row 61
$ws.Range("H61").Value = 2639.4412
$ws.Range("I61").Value = 1747.6
$ws.Range("J61").Value = 3343.5264
$ws.Range("K61").Value = 1747.6
$ws.Range("L61").Value = 3343.5264
$ws.Range("M61").Value = -1535.6
$ws.Range("N61").Value = -3767.5264
# Row 122
$ws.Range("H122").Value = 1545.4736
$ws.Range("I122").Value = 1175
$ws.Range("J122").Value = 1716.4615
$ws.Range("K122").Value = 3525
$ws.Range("L122").Value = 5149.3845
$ws.Range("M122").Value = -1075
$ws.Range("N122").Value = -10049.3845
# Row 136
$ws.Range("H136").Value = 2639.4412
$ws.Range("I136").Value = 1747.6
$ws.Range("J136").Value = 3343.5264
$ws.Range("K136").Value = 5242.799999999999
$ws.Range("L136").Value = 10030.5792
$ws.Range("M136").Value = -2692.799999999999
$ws.Range("N136").Value = -15130.5792

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1527.6923
$ws.Range("I134").Value = 1460
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 4380
$ws.Range("L134").Value = 5700
$ws.Range("M134").Value = -1845
$ws.Range("N134").Value = -10770

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4187.2188
$ws.Range("I31").Value = 3519.5557
$ws.Range("J31").Value = 5045.643
$ws.Range("K31").Value = 3519.5557
$ws.Range("L31").Value = 5045.643
$ws.Range("M31").Value = -3224.5557
$ws.Range("N31").Value = -5635.643
# Row 34
$ws.Range("H34").Value = 4187.2188
$ws.Range("I34").Value = 3519.5557
$ws.Range("J34").Value = 5045.643
$ws.Range("K34").Value = 3519.5557
$ws.Range("L34").Value = 5045.643
$ws.Range("M34").Value = -3317.5557
$ws.Range("N34").Value = -5449.643
# Row 58
$ws.Range("H58").Value = 1767.8864
$ws.Range("I58").Value = 2040.5416
$ws.Range("J58").Value = 1440.7
$ws.Range("K58").Value = 2040.5416
$ws.Range("L58").Value = 1440.7
$ws.Range("M58").Value = -1837.5416
$ws.Range("N58").Value = -1846.7
# Row 99
$ws.Range("H99").Value = 1736.2
$ws.Range("I99").Value = 1408.7273
$ws.Range("K99").Value = 1408.7273
$ws.Range("M99").Value = 89.27269999999999
# Row 126
$ws.Range("H126").Value = 1736.2
$ws.Range("I126").Value = 1408.7273
$ws.Range("K126").Value = 4226.1819
$ws.Range("M126").Value = -1756.1819
# Row 136
$ws.Range("H136").Value = 1767.8864
$ws.Range("I136").Value = 2040.5416
$ws.Range("J136").Value = 1440.7
$ws.Range("K136").Value = 6121.6248
$ws.Range("L136").Value = 4322.1
$ws.Range("M136").Value = -3571.6248
$ws.Range("N136").Value = -9422.1

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 14910.556
$ws.Range("J68").Value = 32650.75
$ws.Range("L68").Value = 97952.25
$ws.Range("N68").Value = -99574.25
# Row 71
$ws.Range("H71").Value = 14910.556
$ws.Range("J71").Value = 32650.75
$ws.Range("L71").Value = 293856.75
$ws.Range("N71").Value = -301968.75
# Row 82
$ws.Range("H82").Value = 12203.75
$ws.Range("J82").Value = 21407.5
$ws.Range("L82").Value = 64222.5
$ws.Range("N82").Value = -65034.5
# Row 85
$ws.Range("H85").Value = 12203.75
$ws.Range("J85").Value = 21407.5
$ws.Range("L85").Value = 64222.5
$ws.Range("N85").Value = -67030.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2091.5
$ws.Range("I80").Value = 1905
$ws.Range("J80").Value = 2171.4285
$ws.Range("K80").Value = 1905
$ws.Range("L80").Value = 2171.4285
$ws.Range("M80").Value = -907
$ws.Range("N80").Value = -4167.4285
# Row 83
$ws.Range("H83").Value = 2091.5
$ws.Range("I83").Value = 1905
$ws.Range("J83").Value = 2171.4285
$ws.Range("K83").Value = 9525
$ws.Range("L83").Value = 10857.1425
$ws.Range("M83").Value = -4533
$ws.Range("N83").Value = -20841.1425
# Row 102
$ws.Range("H102").Value = 2340
$ws.Range("I102").Value = 1914.2858
$ws.Range("J102").Value = 3333.3333
$ws.Range("K102").Value = 1914.2858
$ws.Range("L102").Value = 3333.3333
$ws.Range("M102").Value = -292.2858000000001
$ws.Range("N102").Value = -6577.3333
# Row 113
$ws.Range("H113").Value = 1162.2
$ws.Range("I113").Value = 555.5
$ws.Range("J113").Value = 1566.6666
$ws.Range("K113").Value = 555.5
$ws.Range("L113").Value = 1566.6666
$ws.Range("M113").Value = 1614.5
$ws.Range("N113").Value = -5906.6666
# Row 122
$ws.Range("H122").Value = 1967.2222
$ws.Range("I122").Value = 928.46155
$ws.Range("J122").Value = 2931.7856
$ws.Range("K122").Value = 2785.38465
$ws.Range("L122").Value = 8795.356800000001
$ws.Range("M122").Value = -335.38465
$ws.Range("N122").Value = -13695.3568

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1588.9375
$ws.Range("I7").Value = 1300.5
$ws.Range("J7").Value = 1762
$ws.Range("K7").Value = 1300.5
$ws.Range("L7").Value = 1762
$ws.Range("M7").Value = -1188.5
$ws.Range("N7").Value = -1986
# Row 40
$ws.Range("H40").Value = 3190
$ws.Range("I40").Value = 3093.3333
$ws.Range("J40").Value = 3262.5
$ws.Range("K40").Value = 3093.3333
$ws.Range("L40").Value = 3262.5
$ws.Range("M40").Value = -2957.3333
$ws.Range("N40").Value = -3534.5
# Row 61
$ws.Range("H61").Value = 2760
$ws.Range("I61").Value = 825
$ws.Range("K61").Value = 825
$ws.Range("M61").Value = -623
# Row 113
$ws.Range("H113").Value = 2760
$ws.Range("I113").Value = 825
$ws.Range("K113").Value = 825
$ws.Range("M113").Value = 1345
# Row 126
$ws.Range("H126").Value = 1588.9375
$ws.Range("I126").Value = 1300.5
$ws.Range("J126").Value = 1762
$ws.Range("K126").Value = 3901.5
$ws.Range("L126").Value = 5286
$ws.Range("M126").Value = -1431.5
$ws.Range("N126").Value = -10226
# Row 136
$ws.Range("H136").Value = 3385.8572
$ws.Range("I136").Value = 2975.3333
$ws.Range("J136").Value = 3693.75
$ws.Range("K136").Value = 8925.999899999999
$ws.Range("L136").Value = 11081.25
$ws.Range("M136").Value = -6375.999899999999
$ws.Range("N136").Value = -16181.25

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 286
$ws.Range("I113").Value = 250
$ws.Range("J113").Value = 322
$ws.Range("K113").Value = 750
$ws.Range("L113").Value = 966
$ws.Range("M113").Value = 1420
$ws.Range("N113").Value = -5306
# Row 122
$ws.Range("H122").Value = 323659.72
$ws.Range("I122").Value = 477059.1
$ws.Range("J122").Value = 1521
$ws.Range("K122").Value = 1431177.3
$ws.Range("L122").Value = 4563
$ws.Range("M122").Value = -1428727.3
$ws.Range("N122").Value = -9463
# Row 126
$ws.Range("H126").Value = 417988.22
$ws.Range("I126").Value = 770485.9399999999
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 2311457.82
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -2308987.82
$ws.Range("N126").Value = -9140
